# TimeSheet.xlsx edit: add "Trajectory of points" task block to the
# "مهر 98" (row 48) section, shrinking the task list from 4 rows + totals
# (rows 49-55) down to 3 rows + totals (rows 49-54), matching the sibling
# month block that already lives in rows 38-45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Reshape rows 49-51 (the three task lines) so they carry the same
#    cell styles as the analogous, already-correct block in rows 40-42.
#    Copying a whole cell brings both the style and the (placeholder)
#    value/formula across; we immediately overwrite the value we don't
#    want to keep.
# ---------------------------------------------------------------------

# Row 49: "* BronchoVision GUI", 1 hour, new task note
$ws.Range("B40").Copy($ws.Range("B49"))
$ws.Range("C40").Copy($ws.Range("C49"))
$ws.Range("E40").Copy($ws.Range("E49"))
$ws.Range("C49").Value = 1
$ws.Range("E49").Value = "• Read points from .mat file"

# Row 50: drop the leftover A50 cell entirely (target has no cell there)
$ws.Range("A50").Clear()
$ws.Range("B41").Copy($ws.Range("B50"))
$ws.Range("C41").Copy($ws.Range("C50"))
$ws.Range("E41").Copy($ws.Range("E50"))
$ws.Range("C50").Value = 2
$ws.Range("E50").Value = "• Draw trajectory of points in 3D"

# Row 51: "* Virtual Camera", 4 hours, no note
$ws.Range("B42").Copy($ws.Range("B51"))
$ws.Range("C42").Copy($ws.Range("C51"))
$ws.Range("E42").Copy($ws.Range("E51"))
$ws.Range("C51").Value = 4

# ---------------------------------------------------------------------
# 2) Row 52 becomes the "Total Hours" row (was row 53); its old E52
#    content is gone, and the SUM formula now only spans C49:C51.
# ---------------------------------------------------------------------
$ws.Range("B43").Copy($ws.Range("B52"))
$ws.Range("C43").Copy($ws.Range("C52"))
$ws.Range("E52").Clear()
$ws.Range("C52").Formula = "=SUM(C49:C51)"

# ---------------------------------------------------------------------
# 3) Rows 53 & 54 become the "Paid"/"Not Paid" summary (previously rows
#    54 & 55); drop the stale B53 task-label cell.
# ---------------------------------------------------------------------
$ws.Range("B53").Clear()
$ws.Range("C54").Copy($ws.Range("C53"))
$ws.Range("D54").Copy($ws.Range("D53"))
$ws.Range("C55").Copy($ws.Range("C54"))
$ws.Range("D55").Copy($ws.Range("D54"))
$ws.Range("D53").Value = 3
$ws.Range("D54").Value = 4

# Row 55 no longer exists in the sheet - remove it entirely.
$ws.Rows.Item(55).Delete() | Out-Null

# ---------------------------------------------------------------------
# 4) Update the view so the newly edited block is in frame, matching
#    the author's on-save cursor position.
# ---------------------------------------------------------------------
$ws.Range("E52").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1
